# Convert HOUR_APPR_PROCESS_START (column V) values from plain numbers
# into text strings formatted as "HH:00:00" (e.g. 19 -> "19:00:00").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column V is column 22 ("HOUR_APPR_PROCESS_START"); data rows are 2..6.
for ($row = 2; $row -le 6; $row++) {
    $cell = $ws.Cells.Item($row, 22)
    $hour = $cell.Value()
    $text = [string]([int]$hour) + ":00:00"
    $cell.Value = $text
}
